$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 171.934662575
$ws.Range("C2").Value = 20.8375089279

$ws.Range("B3").Value = 171934.662575
$ws.Range("C3").Value = 20837.5089279

$ws.Range("B4").Value = 343869.32515
$ws.Range("C4").Value = 62512.5267837

$ws.Range("B5").Value = 6877.386503
$ws.Range("C5").Value = 1250.250535674
